# Update Excel file with latest predictions
$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Home win": refresh the 4 existing prediction rows with new matches
# ---------------------------------------------------------------------------
$wsHome = $wb.Worksheets.Item("Home win")

$wsHome.Range("A2").Value = "14-01-2025 13:00"
$wsHome.Range("B2").Value = "ALGERIA"
$wsHome.Range("C2").Value = "COUPE NATIONALE"
$wsHome.Range("D2").Value = "RC Kouba - JS Saoura"
$wsHome.Range("E2").Value = 70
$wsHome.Range("F2").Value = 2.8

$wsHome.Range("A3").Value = "14-01-2025 19:45"
$wsHome.Range("B3").Value = "ENGLAND"
$wsHome.Range("C3").Value = "NON LEAGUE PREMIER - SOUTHERN CENTRAL"
$wsHome.Range("D3").Value = "Harborough Town - AFC Sudbury"
$wsHome.Range("E3").Value = 80
$wsHome.Range("F3").Value = 1.95

$wsHome.Range("A4").Value = "15-01-2025 18:30"
$wsHome.Range("B4").Value = "SPAIN"
$wsHome.Range("C4").Value = "COPA DEL REY"
$wsHome.Range("D4").Value = "Pontevedra - Getafe"
$wsHome.Range("E4").Value = 70
$wsHome.Range("F4").Value = 3.5

$wsHome.Range("A5").Value = "15-01-2025 12:00"
$wsHome.Range("B5").Value = "THAILAND"
$wsHome.Range("C5").Value = "THAI LEAGUE 1"
$wsHome.Range("D5").Value = "Bangkok Glass - Ratchaburi"
$wsHome.Range("E5").Value = 73.3
$wsHome.Range("F5").Value = 1.7

# ---------------------------------------------------------------------------
# Sheet "Draw": row 2 stays, row 3 is refreshed, 3 new rows are appended
# ---------------------------------------------------------------------------
$wsDraw = $wb.Worksheets.Item("Draw")

$wsDraw.Range("A3").Value = "14-01-2025 23:00"
$wsDraw.Range("B3").Value = "BRAZIL"
$wsDraw.Range("C3").Value = "SERGIPANO"
$wsDraw.Range("D3").Value = "Barra SE - Lagarto"
$wsDraw.Range("E3").Value = 60
$wsDraw.Range("F3").Value = 4.5

$wsDraw.Range("A4").Value = "14-01-2025 12:30"
$wsDraw.Range("B4").Value = "EGYPT"
$wsDraw.Range("C4").Value = "SECOND LEAGUE"
$wsDraw.Range("D4").Value = "Suez - Abu Qair Semad"
$wsDraw.Range("E4").Value = 70
$wsDraw.Range("F4").Value = 2.88

$wsDraw.Range("A5").Value = "14-01-2025 12:00"
$wsDraw.Range("B5").Value = "ETHIOPIA"
$wsDraw.Range("C5").Value = "PREMIER LEAGUE"
$wsDraw.Range("D5").Value = "Ethiopia Nigd Bank - Mebrat Hayl"
$wsDraw.Range("E5").Value = 66.7
$wsDraw.Range("F5").Value = 2.88

$wsDraw.Range("A6").Value = "14-01-2025 11:00"
$wsDraw.Range("B6").Value = "THAILAND"
$wsDraw.Range("C6").Value = "THAI LEAGUE 1"
$wsDraw.Range("D6").Value = "Lamphun Warrior - Nakhon Ratchasima FC"
$wsDraw.Range("E6").Value = 60
$wsDraw.Range("F6").Value = 3.75

# ---------------------------------------------------------------------------
# Sheet "Btts": rows 2-5 refreshed, 5 new rows appended
# ---------------------------------------------------------------------------
$wsBtts = $wb.Worksheets.Item("Btts")

$wsBtts.Range("A2").Value = "14-01-2025 19:45"
$wsBtts.Range("B2").Value = "ITALY"
$wsBtts.Range("C2").Value = "SERIE A"
$wsBtts.Range("D2").Value = "Atalanta - Juventus"
$wsBtts.Range("E2").Value = 76.7
$wsBtts.Range("F2").Value = 1.8

$wsBtts.Range("A3").Value = "14-01-2025 19:45"
$wsBtts.Range("B3").Value = "ENGLAND"
$wsBtts.Range("C3").Value = "NATIONAL LEAGUE - NORTH"
$wsBtts.Range("D3").Value = "Chorley - Peterborough Sports"
$wsBtts.Range("E3").Value = 76.7
$wsBtts.Range("F3").Value = 1.75

$wsBtts.Range("A4").Value = "14-01-2025 19:45"
$wsBtts.Range("B4").Value = "FRANCE"
$wsBtts.Range("C4").Value = "COUPE DE FRANCE"
$wsBtts.Range("D4").Value = "Haguenau - Dunkerque"
$wsBtts.Range("E4").Value = 80
$wsBtts.Range("F4").Value = 1.93

$wsBtts.Range("A5").Value = "15-01-2025 19:30"
$wsBtts.Range("B5").Value = "GERMANY"
$wsBtts.Range("C5").Value = "BUNDESLIGA"
$wsBtts.Range("D5").Value = "Union Berlin - FC Augsburg"
$wsBtts.Range("E5").Value = 80
$wsBtts.Range("F5").Value = 1.95

$wsBtts.Range("A6").Value = "15-01-2025 21:30"
$wsBtts.Range("B6").Value = "BRAZIL"
$wsBtts.Range("C6").Value = "PAULISTA - A1"
$wsBtts.Range("D6").Value = "Velo Clube - Noroeste"
$wsBtts.Range("E6").Value = 81.7
$wsBtts.Range("F6").Value = 2.05

$wsBtts.Range("A7").Value = "15-01-2025 23:00"
$wsBtts.Range("B7").Value = "COSTA-RICA"
$wsBtts.Range("C7").Value = "PRIMERA DIVISIÓN"
$wsBtts.Range("D7").Value = "Santa Ana - Puntarenas FC"
$wsBtts.Range("E7").Value = 83.3
$wsBtts.Range("F7").Value = 1.83

$wsBtts.Range("A8").Value = "15-01-2025 17:30"
$wsBtts.Range("B8").Value = "FRANCE"
$wsBtts.Range("C8").Value = "COUPE DE FRANCE"
$wsBtts.Range("D8").Value = "Thaon - Strasbourg"
$wsBtts.Range("E8").Value = 84
$wsBtts.Range("F8").Value = 2.6

$wsBtts.Range("A9").Value = "15-01-2025 17:00"
$wsBtts.Range("B9").Value = "SAUDI-ARABIA"
$wsBtts.Range("C9").Value = "PRO LEAGUE"
$wsBtts.Range("D9").Value = "Al Kholood - Al-Ahli Jeddah"
$wsBtts.Range("E9").Value = 76.7
$wsBtts.Range("F9").Value = 1.85

$wsBtts.Range("A10").Value = "15-01-2025 11:00"
$wsBtts.Range("B10").Value = "THAILAND"
$wsBtts.Range("C10").Value = "THAI LEAGUE 1"
$wsBtts.Range("D10").Value = "Port FC - Khon Kaen United"
$wsBtts.Range("E10").Value = 87.8
$wsBtts.Range("F10").Value = 1.75

# ---------------------------------------------------------------------------
# Sheet "Over_Under": rows 2-5 untouched, 2 new rows appended
# ---------------------------------------------------------------------------
$wsOU = $wb.Worksheets.Item("Over_Under")

$wsOU.Range("A6").Value = "15-01-2025 17:45"
$wsOU.Range("B6").Value = "NETHERLANDS"
$wsOU.Range("C6").Value = "KNVB BEKER"
$wsOU.Range("D6").Value = "GO Ahead Eagles - Twente"
$wsOU.Range("E6").Value = 80
$wsOU.Range("F6").Value = 1.73
$wsOU.Range("G6").Value = 50
$wsOU.Range("H6").Value = 2.8

$wsOU.Range("A7").Value = "15-01-2025 17:30"
$wsOU.Range("B7").Value = "FRANCE"
$wsOU.Range("C7").Value = "COUPE DE FRANCE"
$wsOU.Range("D7").Value = "Cannes - Lorient"
$wsOU.Range("E7").Value = 80
$wsOU.Range("F7").Value = 1.8
$wsOU.Range("G7").Value = 60
$wsOU.Range("H7").Value = 3

# "Away Win" sheet is unchanged in this update.
